$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I6").Value = 'sd'
$ws.Range("J6").Value = 'Statement-non-opinion'
$ws.Range("I7").Value = 'sv'
$ws.Range("J7").Value = 'Statement-opinion'
$ws.Range("I10").Value = 'sd'
$ws.Range("J10").Value = 'Statement-non-opinion'
$ws.Range("I14").Value = 'b'
$ws.Range("J14").Value = 'Acknowledge (Backchannel)'
$ws.Range("I22").Value = 'sd'
$ws.Range("J22").Value = 'Statement-non-opinion'
$ws.Range("I23").Value = 'sd'
$ws.Range("J23").Value = 'Statement-non-opinion'
$ws.Range("I25").Value = 'sd'
$ws.Range("J25").Value = 'Statement-non-opinion'
$ws.Range("I52").Value = 'sd'
$ws.Range("J52").Value = 'Statement-non-opinion'
$ws.Range("I61").Value = 'ba'
$ws.Range("J61").Value = 'Appreciation'
$ws.Range("I63").Value = 'sd'
$ws.Range("J63").Value = 'Statement-non-opinion'
$ws.Range("I64").Value = 'aa'
$ws.Range("J64").Value = 'Agree/Accept'
$ws.Range("I70").Value = 'sv'
$ws.Range("J70").Value = 'Statement-opinion'
$ws.Range("I80").Value = 'aa'
$ws.Range("J80").Value = 'Agree/Accept'
$ws.Range("I86").Value = 'b'
$ws.Range("J86").Value = 'Acknowledge (Backchannel)'
$ws.Range("I113").Value = 'sd'
$ws.Range("J113").Value = 'Statement-non-opinion'
$ws.Range("I124").Value = 'sd'
$ws.Range("J124").Value = 'Statement-non-opinion'
$ws.Range("I126").Value = 'ba'
$ws.Range("J126").Value = 'Appreciation'
$ws.Range("I129").Value = 'sd'
$ws.Range("J129").Value = 'Statement-non-opinion'
$ws.Range("I134").Value = 'b'
$ws.Range("J134").Value = 'Acknowledge (Backchannel)'
$ws.Range("I151").Value = 'aa'
$ws.Range("J151").Value = 'Agree/Accept'
$ws.Range("I166").Value = 'sv'
$ws.Range("J166").Value = 'Statement-opinion'
$ws.Range("I180").Value = 'sv'
$ws.Range("J180").Value = 'Statement-opinion'
$ws.Range("I187").Value = 'sd'
$ws.Range("J187").Value = 'Statement-non-opinion'
$ws.Range("I188").Value = 'sv'
$ws.Range("J188").Value = 'Statement-opinion'
$ws.Range("I190").Value = 'aa'
$ws.Range("J190").Value = 'Agree/Accept'
$ws.Range("I197").Value = 'aa'
$ws.Range("J197").Value = 'Agree/Accept'
$ws.Range("I209").Value = 'b'
$ws.Range("J209").Value = 'Acknowledge (Backchannel)'
$ws.Range("I213").Value = 'sd'
$ws.Range("J213").Value = 'Statement-non-opinion'
$ws.Range("I221").Value = '%'
$ws.Range("J221").Value = 'Uninterpretable'
$ws.Range("I248").Value = 'sd'
$ws.Range("J248").Value = 'Statement-non-opinion'
$ws.Range("I263").Value = 'sd'
$ws.Range("J263").Value = 'Statement-non-opinion'
$ws.Range("I303").Value = 'aa'
$ws.Range("J303").Value = 'Agree/Accept'
$ws.Range("I304").Value = 'b'
$ws.Range("J304").Value = 'Acknowledge (Backchannel)'
$ws.Range("I311").Value = 'b'
$ws.Range("J311").Value = 'Acknowledge (Backchannel)'
$ws.Range("I329").Value = 'sd'
$ws.Range("J329").Value = 'Statement-non-opinion'
$ws.Range("I336").Value = 'sv'
$ws.Range("J336").Value = 'Statement-opinion'
$ws.Range("I346").Value = 'aa'
$ws.Range("J346").Value = 'Agree/Accept'
$ws.Range("I348").Value = 'sv'
$ws.Range("J348").Value = 'Statement-opinion'
$ws.Range("I357").Value = 'sd'
$ws.Range("J357").Value = 'Statement-non-opinion'
$ws.Range("I382").Value = 'aa'
$ws.Range("J382").Value = 'Agree/Accept'
$ws.Range("I384").Value = 'aa'
$ws.Range("J384").Value = 'Agree/Accept'
$ws.Range("I386").Value = 'sd'
$ws.Range("J386").Value = 'Statement-non-opinion'
$ws.Range("I388").Value = 'ba'
$ws.Range("J388").Value = 'Appreciation'
$ws.Range("I391").Value = 'sv'
$ws.Range("J391").Value = 'Statement-opinion'
$ws.Range("I396").Value = 'sd'
$ws.Range("J396").Value = 'Statement-non-opinion'
$ws.Range("I398").Value = 'aa'
$ws.Range("J398").Value = 'Agree/Accept'
$ws.Range("I404").Value = 'sv'
$ws.Range("J404").Value = 'Statement-opinion'
$ws.Range("I407").Value = 'sd'
$ws.Range("J407").Value = 'Statement-non-opinion'
$ws.Range("I408").Value = 'sd'
$ws.Range("J408").Value = 'Statement-non-opinion'
$ws.Range("I418").Value = 'sv'
$ws.Range("J418").Value = 'Statement-opinion'
$ws.Range("I427").Value = 'sv'
$ws.Range("J427").Value = 'Statement-opinion'
$ws.Range("I429").Value = 'aa'
$ws.Range("J429").Value = 'Agree/Accept'
$ws.Range("I437").Value = 'aa'
$ws.Range("J437").Value = 'Agree/Accept'
$ws.Range("I448").Value = 'sd'
$ws.Range("J448").Value = 'Statement-non-opinion'
$ws.Range("I465").Value = 'sd'
$ws.Range("J465").Value = 'Statement-non-opinion'
$ws.Range("I468").Value = 'sv'
$ws.Range("J468").Value = 'Statement-opinion'
$ws.Range("I470").Value = 'sd'
$ws.Range("J470").Value = 'Statement-non-opinion'
$ws.Range("I471").Value = 'aa'
$ws.Range("J471").Value = 'Agree/Accept'
$ws.Range("I474").Value = 'ba'
$ws.Range("J474").Value = 'Appreciation'
$ws.Range("I476").Value = 'sd'
$ws.Range("J476").Value = 'Statement-non-opinion'
$ws.Range("I480").Value = '%'
$ws.Range("J480").Value = 'Uninterpretable'
$ws.Range("I489").Value = 'aa'
$ws.Range("J489").Value = 'Agree/Accept'
$ws.Range("I490").Value = 'sd'
$ws.Range("J490").Value = 'Statement-non-opinion'
$ws.Range("I491").Value = 'b'
$ws.Range("J491").Value = 'Acknowledge (Backchannel)'
$ws.Range("I495").Value = 'sd'
$ws.Range("J495").Value = 'Statement-non-opinion'
$ws.Range("I502").Value = 'sd'
$ws.Range("J502").Value = 'Statement-non-opinion'
$ws.Range("I506").Value = 'aa'
$ws.Range("J506").Value = 'Agree/Accept'
$ws.Range("I524").Value = '%'
$ws.Range("J524").Value = 'Uninterpretable'
$ws.Range("I533").Value = 'sd'
$ws.Range("J533").Value = 'Statement-non-opinion'
$ws.Range("I539").Value = 'sv'
$ws.Range("J539").Value = 'Statement-opinion'
$ws.Range("I552").Value = 'ba'
$ws.Range("J552").Value = 'Appreciation'
$ws.Range("I555").Value = 'sd'
$ws.Range("J555").Value = 'Statement-non-opinion'
